$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''56.089.21'
$ws.Range('E2').Value = '  +9.19%  '
$ws.Range('D3').Value = '''3.220.37'
$ws.Range('E3').Value = '  +4.16%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''396.53'
$ws.Range('E5').Value = '  +2.57%  '
$ws.Range('D6').Value = '''110.86'
$ws.Range('E6').Value = '  +7.11%  '
$ws.Range('E7').Value = '  +2.57%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +5.35%  '
$ws.Range('D10').Value = '''39.17'
$ws.Range('E10').Value = '  +6.15%  '
$ws.Range('D11').Value = '''0.0914'
$ws.Range('E11').Value = '  +6.82%  '
$ws.Range('E12').Value = '  +2.03%  '
$ws.Range('D13').Value = '''3.724.86'
$ws.Range('E13').Value = '  +4.14%  '
$ws.Range('E14').Value = '  +3.70%  '
$ws.Range('D15').Value = '''19.01'
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('D16').Value = '''3.190.56'
$ws.Range('E16').Value = '  +3.02%  '
$ws.Range('D17').Value = '''1.04'
$ws.Range('E17').Value = '  +5.29%  '
$ws.Range('D18').Value = '''10.83'
$ws.Range('E18').Value = '  +1.76%  '
$ws.Range('D19').Value = '''55.895.31'
$ws.Range('E19').Value = '  +8.66%  '
$ws.Range('D20').Value = '''3.33'
$ws.Range('E20').Value = '  +3.41%  '
$ws.Range('D21').Value = '''0.0000103'
$ws.Range('E21').Value = '  +6.55%  '
$ws.Range('D22').Value = '''12.95'
$ws.Range('E22').Value = '  +3.78%  '
$ws.Range('D23').Value = '''299.05'
$ws.Range('E23').Value = '  +12.63%  '
$ws.Range('D24').Value = '''75.46'
$ws.Range('E24').Value = '  +7.85%  '
$ws.Range('D25').Value = '''3.21'
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('D26').Value = '''8.15'
$ws.Range('E26').Value = '  +1.93%  '
$ws.Range('D27').Value = '''28.12'
$ws.Range('E27').Value = '  +2.87%  '
$ws.Range('D28').Value = '''7.52'
$ws.Range('E28').Value = '  +4.50%  '
$ws.Range('E29').Value = '  +4.58%  '
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('E31').Value = '  +3.95%  '
$ws.Range('D32').Value = '''11.12'
$ws.Range('E32').Value = '  +6.88%  '
$ws.Range('E33').Value = '  +3.44%  '
$ws.Range('D34').Value = '''36.13'
$ws.Range('E34').Value = '  +1.41%  '
$ws.Range('D35').Value = '''2.12'
$ws.Range('E35').Value = '  +2.48%  '
$ws.Range('E36').Value = '  +2.78%  '
$ws.Range('D37').Value = '''3.13'
$ws.Range('E37').Value = '  +25.87%  '
$ws.Range('D38').Value = '''3.53'
$ws.Range('E38').Value = '  +5.07%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '''135.59'
$ws.Range('E40').Value = '  +4.70%  '
$ws.Range('D41').Value = '''17.35'
$ws.Range('E41').Value = '  +5.15%  '
$ws.Range('E42').Value = '  +3.18%  '
$ws.Range('D43').Value = '''3.99'
$ws.Range('E43').Value = '  +4.35%  '
$ws.Range('E44').Value = '  +3.02%  '
$ws.Range('D45').Value = '''0.283'
$ws.Range('E45').Value = '  -2.89%  '
$ws.Range('D46').Value = '''22.15'
$ws.Range('E46').Value = '  +0.51%  '
$ws.Range('E47').Value = '  +52.48%  '
$ws.Range('E48').Value = '  +1.53%  '
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('D50').Value = '''2.119.02'
$ws.Range('E50').Value = '  +2.21%  '
$ws.Range('D51').Value = '''0.0359'
$ws.Range('E51').Value = '  +8.86%  '
